# stocks_bought.xlsx: correct the DOCU row's buy_fees_eur (D8)
# from 1.65 to 42.39.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("D8").Value = 42.39

# Leave the selection where the editor last left it.
$ws.Range("G8").Select()
